# VM_instructions.pptx - slide 11, "TextBox 2" (Step 6 instructions)
# Split the run that read:
#   " notebook server is running on your VM you can then open a browser and
#    copy and paste the URL that starts with "
# into three runs (same Arial/en-AU formatting) reading:
#   " notebook server is running on your VM you can then open a browser and "
#   "copy then "
#   "paste the URL that starts with "

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$full = $tr.Text
$oldRun = " notebook server is running on your VM you can then open a browser and copy and paste the URL that starts with "
$idx = $full.IndexOf($oldRun)

if ($idx -ge 0) {
    $partA_old = " notebook server is running on your VM you can then open a browser and "
    $partB_old = "copy and "
    $partC_old = "paste the URL that starts with "

    $partA_new = " notebook server is running on your VM you can then open a browser and "
    $partB_new = "copy then "
    $partC_new = "paste the URL that starts with "

    $startA = $idx + 1
    $startB = $startA + $partA_old.Length
    $startC = $startB + $partB_old.Length

    # Apply edits right-to-left so offsets computed from the original text
    # stay valid (partB's length changes from 9 to 10 characters).
    $subC = $tr.Characters($startC, $partC_old.Length)
    $subC.Text = $partC_new

    $subB = $tr.Characters($startB, $partB_old.Length)
    $subB.Text = $partB_new

    $subA = $tr.Characters($startA, $partA_old.Length)
    $subA.Text = $partA_new
}
